$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers (F1, G1) and data (F2:G4) for the new
# "Item disc" / "Category Disc" columns.
$ws.Range("F1").Value = "Item disc"
$ws.Range("G1").Value = "Category Disc"
$ws.Range("G1").HorizontalAlignment = -4108

$ws.Range("F2").Value = 10
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("G2").Value = 10
$ws.Range("G2").HorizontalAlignment = -4108

$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 15
$ws.Range("G3").HorizontalAlignment = -4108

$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 20

# Update the selection/view: select G7 (this also drops the old
# topLeftCell scroll position that pointed at A7).
$ws.Range("G7").Select()
